$d = $word.ActiveDocument

# Find the "Author" paragraph containing "Edison Achalma" (the one right
# after the "Editar: Editar" Heading1, not the later mentions in the
# Nota de Autores / CRediT / correspondence sections).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Edison Achalma" -and `
        $p.Style.NameLocal -eq "Author") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $r = $target.Range
    $r.InsertAfter([char]13 + "Escuela Profesional de Economía, Universidad Nacional de San Cristóbal de Huamanga")
}
